$d = $word.ActiveDocument

# Update the title/date paragraph (first paragraph, outside the table)
$titleRange = $d.Paragraphs.Item(1).Range
$titleRange.Text = "2023-11-22 Wednesday"

# Update each table cell in row-major order to match document order
$t = $d.Tables.Item(1)
$newValues = @(
    "39+37=",
    "46+3=",
    "4-2=",
    "93-63=",
    "15-11=",
    "65-49=",
    "36-15=",
    "73-50=",
    "67+3=",
    "20+71=",
    "31+23=",
    "4+1=",
    "41+5=",
    "80-44=",
    "15+63=",
    "64-19=",
    "49-3=",
    "71+10=",
    "63-0=",
    "54+10=",
    "41+42=",
    "89-52=",
    "90-1=",
    "57+36=",
    "34+44=",
    "51+23=",
    "34-21=",
    "0+16=",
    "24+41=",
    "83-47=",
    "75+5=",
    "74-24=",
    "92-57=",
    "83-67=",
    "28-3=",
    "98-48=",
    "44-11=",
    "84+5=",
    "18-5=",
    "5+14=",
    "37-3=",
    "22+45=",
    "44+9=",
    "54+6=",
    "53+31=",
    "71+27=",
    "56+23=",
    "15+84=",
    "63-38=",
    "71-1=",
    "13+67=",
    "63-57=",
    "21+61=",
    "60-37=",
    "6+70=",
    "63+35=",
    "39+27=",
    "20-14=",
    "71-20=",
    "34+45=",
    "48+0=",
    "21+62=",
    "88+6=",
    "0+54=",
    "1+11=",
    "65-26=",
    "57-47=",
    "67-15=",
    "88-59=",
    "16+31=",
    "10+71=",
    "71-53=",
    "99-42=",
    "89+2=",
    "14+2=",
    "80-64=",
    "68-7=",
    "56-47=",
    "16+75=",
    "57-22=",
    "74+9=",
    "66+15=",
    "3+70=",
    "8+16=",
    "83-63=",
    "34+5=",
    "8+24=",
    "56-38=",
    "49+22=",
    "86-6=",
    "87-73=",
    "65-9=",
    "43-36=",
    "83-20=",
    "14+48=",
    "58-1=",
    "35-28=",
    "40+2=",
    "99-38=",
    "88-1="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output "Done: updated $idx cells"
